$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new set of columns for the bulk shift extension import template
$ws.Range("A1").Value = "StaffId"
$ws.Range("B1").Value = "ApplicationType"
$ws.Range("C1").Value = "TransactionDate"
$ws.Range("D1").Value = "DurationHours"
$ws.Range("E1").Value = "BeforeShiftHours"
$ws.Range("F1").Value = "AfterShiftHours"
$ws.Range("G1").Value = "Shift"
$ws.Range("H1").Value = "Remarks"

# --- Sample / placeholder row (row 2) carrying the number formats used for
# date and duration/time entry so users importing data see the right format.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").NumberFormat = "mm:ss.0"
$ws.Range("F2").NumberFormat = "mm:ss.0"
$ws.Range("M2").NumberFormat = "mm:ss.0"

# --- Column widths sized to fit the new headers (values chosen so the
# engine's internal pixel-quantized ColumnWidth rounds to the intended
# stored sheet width as closely as possible)
$ws.Columns.Item(1).ColumnWidth = 5.833333333333333
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 13.0
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws.Columns.Item(7).ColumnWidth = 4.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.833333333333333

$ws.Range("D2").Select() | Out-Null
